# Add two new rows of daily OHLCV data (27 Feb 2020 and 28 Feb 2020) to the
# bottom of the existing "ISTONE" (ticker 0209) price history table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 145; Timestamp = 1582761600; Date = "2020-02-27"; Id = "0209"; Name = "ISTONE"; Open = 0.225; High = 0.235; Low = 0.22;  Close = 0.22;  Vol = 17753100 },
    @{ Row = 146; Timestamp = 1582848000; Date = "2020-02-28"; Id = "0209"; Name = "ISTONE"; Open = 0.215; High = 0.215; Low = 0.195; Close = 0.195; Vol = 21701900 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Numeric columns: timestamp, open, high, low, close, volume.
    $ws.Range("A$row").Value = $r.Timestamp
    $ws.Range("E$row").Value = $r.Open
    $ws.Range("F$row").Value = $r.High
    $ws.Range("G$row").Value = $r.Low
    $ws.Range("H$row").Value = $r.Close
    $ws.Range("I$row").Value = $r.Vol

    # Text columns: force text format first so values like "2020-02-27" and
    # "0209" are stored as literal strings instead of being auto-converted
    # to a date serial / number, then restore the default "Normal" style so
    # no stray cell formatting is introduced.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r.Date
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = $r.Id
    $ws.Range("C$row").Style = "Normal"

    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $r.Name
    $ws.Range("D$row").Style = "Normal"
}
